$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest representable value given pixel-grid quantization
# of the ColumnWidth property; target stored width is 14.42578125)
$ws.Columns.Item(1).ColumnWidth = 13.67
$ws.Columns.Item(2).ColumnWidth = 13.67

# Update cell values
$ws.Range("A1").Value = -0.036093149507682783
$ws.Range("B1").Value = 0.036093148796240378

$ws.Range("A2").Value = 0.015304657173543283
$ws.Range("B2").Value = -0.015304657916298717

$ws.Range("A3").Value = 0.014685455106506441
$ws.Range("B3").Value = -0.014685455820034303

$ws.Range("A4").Value = -0.053678655900885358
$ws.Range("B4").Value = 0.053678655182271995
